$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header labels so the separator changes from "." to ",".
$headers = @(
    "Team",
    "Points",
    "Matches",
    "Wins",
    "Draws",
    "Loses",
    "Goals,scored",
    "Goals,conceded",
    "Difference,goals",
    "Percentage,scored,goals",
    "Percentage,conceded,goals",
    "Shots",
    "Shots,on,goal",
    "Penalties,scored",
    "Assistances",
    "Fouls,made",
    "Matches,without,conceding",
    "Yellow,cards",
    "Red,cards",
    "Offsides"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Columns J (Percentage.scored.goals) and K (Percentage.conceded.goals) were
# stored as text; turn them into real numbers for every data row.
$percentages = @(
    @(2.95, 0.76),
    @(2.89, 0.89),
    @(1.66, 0.47),
    @(1.16, 0.92),
    @(1.53, 1.18),
    @(1.34, 1.55),
    @(1.34, 1.32),
    @(1, 0.92),
    @(1.18, 1.26),
    @(0.89, 1.37),
    @(1.18, 1.39),
    @(1.21, 1.26),
    @(1.29, 1.61),
    @(1.05, 1.95),
    @(1.18, 1.61),
    @(1.21, 1.82),
    @(1.05, 1.63),
    @(1.37, 1.92),
    @(0.97, 1.76),
    @(0.97, 1.84)
)
for ($i = 0; $i -lt $percentages.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $percentages[$i][0]
    $ws.Cells.Item($row, 11).Value = $percentages[$i][1]
}

# Move the active selection to G12, matching the saved view state.
$null = $ws.Range("G12").Select()
